$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Paul Pot"
$ws.Range("C4").Value = "po@yandex.ru"
$ws.Range("D4").Value = "23cst4"
